$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 2, shifting all existing
# data rows (2-11) down to (3-12).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new observation record (Jarpe / Tetrastes
# bonasia, id 108660302).
$ws.Cells.Item(2, 1).Value = 108660302
$ws.Cells.Item(2, 2).Value = 57064
$ws.Cells.Item(2, 4).Value = "NT"
$ws.Cells.Item(2, 5).Value = 102612
$ws.Cells.Item(2, 6).Value = "Järpe"
$ws.Cells.Item(2, 7).Value = "Tetrastes bonasia"
$ws.Cells.Item(2, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(2, 13).Value = "äldre spillning"
$ws.Cells.Item(2, 16).Value = "Marsbäcken, väster Kullavägen, Vb"
$ws.Cells.Item(2, 17).Value = 752587
$ws.Cells.Item(2, 18).Value = 7093727
$ws.Cells.Item(2, 19).Value = 25
$ws.Cells.Item(2, 20).Value = "Västerbotten"
$ws.Cells.Item(2, 21).Value = "Umeå"
$ws.Cells.Item(2, 22).Value = "Västerbotten"
$ws.Cells.Item(2, 23).Value = "Umeå socken"
$ws.Cells.Item(2, 25).Value = "'2023-05-01"
$ws.Cells.Item(2, 26).Value = "13:00"
$ws.Cells.Item(2, 27).Value = "'2023-05-01"
$ws.Cells.Item(2, 28).Value = "15:45"
$ws.Cells.Item(2, 29).Value = "Hittades på två olika ställen"
$ws.Cells.Item(2, 30).Value = $false
$ws.Cells.Item(2, 31).Value = $false
$ws.Cells.Item(2, 33).Value = $false
$ws.Cells.Item(2, 49).Value = "Roger Olofsson"
$ws.Cells.Item(2, 50).Value = "Roger Olofsson"
